$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4867.3115
$ws.Range("I138").Value = 1691.8
$ws.Range("J138").Value = 7940.387
$ws.Range("K138").Value = 5075.4
$ws.Range("L138").Value = 23821.161
$ws.Range("M138").Value = 64.60000000000036
$ws.Range("N138").Value = -34101.161
$ws.Range("H141").Value = 3819.5625
$ws.Range("I141").Value = 1434
$ws.Range("J141").Value = 10976.25
$ws.Range("K141").Value = 4302
$ws.Range("L141").Value = 32928.75
$ws.Range("M141").Value = 878
$ws.Range("N141").Value = -43288.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6066960
$ws.Range("I32").Value = 5862.9067
$ws.Range("J32").Value = 27785892
$ws.Range("K32").Value = 5862.9067
$ws.Range("L32").Value = 27785892
$ws.Range("M32").Value = -5575.9067
$ws.Range("N32").Value = -27786466
$ws.Range("H132").Value = 1091561.5
$ws.Range("I132").Value = 2065.275
$ws.Range("J132").Value = 4204408
$ws.Range("K132").Value = 6195.825000000001
$ws.Range("L132").Value = 12613224
$ws.Range("M132").Value = -3665.825000000001
$ws.Range("N132").Value = -12618284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1242.7317
$ws.Range("I31").Value = 778.70966
$ws.Range("J31").Value = 1524.7843
$ws.Range("K31").Value = 778.70966
$ws.Range("L31").Value = 1524.7843
$ws.Range("M31").Value = -483.70966
$ws.Range("N31").Value = -2114.7843
$ws.Range("H34").Value = 1242.7317
$ws.Range("I34").Value = 778.70966
$ws.Range("J34").Value = 1524.7843
$ws.Range("K34").Value = 778.70966
$ws.Range("L34").Value = 1524.7843
$ws.Range("M34").Value = -576.70966
$ws.Range("N34").Value = -1928.7843
$ws.Range("H134").Value = 13890025
$ws.Range("I134").Value = 1060.24
$ws.Range("J134").Value = 45455856
$ws.Range("K134").Value = 3180.72
$ws.Range("L134").Value = 136367568
$ws.Range("M134").Value = -645.7200000000003
$ws.Range("N134").Value = -136372638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5553.9375
$ws.Range("I68").Value = 1333.3334
$ws.Range("J68").Value = 5990.552
$ws.Range("K68").Value = 4000.0002
$ws.Range("L68").Value = 17971.656
$ws.Range("M68").Value = -3189.0002
$ws.Range("N68").Value = -19593.656
$ws.Range("H71").Value = 5553.9375
$ws.Range("I71").Value = 1333.3334
$ws.Range("J71").Value = 5990.552
$ws.Range("K71").Value = 12000.0006
$ws.Range("L71").Value = 53914.96799999999
$ws.Range("M71").Value = -7944.000599999999
$ws.Range("N71").Value = -62026.96799999999
$ws.Range("H113").Value = 5193841.5
$ws.Range("I113").Value = 3205824
$ws.Range("J113").Value = 7778264
$ws.Range("K113").Value = 9617472
$ws.Range("L113").Value = 23334792
$ws.Range("M113").Value = -9615302
$ws.Range("N113").Value = -23339132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H126").Value = 7145314.5
$ws.Range("I126").Value = 1970.6666
$ws.Range("J126").Value = 9093499
$ws.Range("K126").Value = 5911.9998
$ws.Range("L126").Value = 27280497
$ws.Range("M126").Value = -3441.9998
$ws.Range("N126").Value = -27285437

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119:N119").ClearContents()
$ws.Range("H120:N120").ClearContents()
$ws.Range("H121:N121").ClearContents()
$ws.Range("H122:N122").ClearContents()
$ws.Range("H123:N123").ClearContents()
$ws.Range("H124:N124").ClearContents()
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H135:N135").ClearContents()
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()
